$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 295.875  # H39: 199 -> 295.875
$ws.Cells.Item(39, 9).Value = 361.16666  # I39: 218.8 -> 361.16666
$ws.Cells.Item(39, 11).Value = 1083.49998  # K39: 656.4000000000001 -> 1083.49998
$ws.Cells.Item(39, 13).Value = -787.4999800000001  # M39: -360.4000000000001 -> -787.4999800000001
$ws.Cells.Item(69, 8).Value = 6000  # H69: 0 -> 6000
$ws.Cells.Item(69, 10).Value = 6000  # J69: 0 -> 6000
$ws.Cells.Item(69, 12).Value = 18000  # L69: 0 -> 18000
$ws.Cells.Item(69, 14).Value = -19748  # N69: None -> -19748
$ws.Cells.Item(72, 8).Value = 6000  # H72: 0 -> 6000
$ws.Cells.Item(72, 10).Value = 6000  # J72: 0 -> 6000
$ws.Cells.Item(72, 12).Value = 54000  # L72: 0 -> 54000
$ws.Cells.Item(72, 14).Value = -62736  # N72: None -> -62736
$ws.Cells.Item(107, 8).Value = 467.6316  # H107: 368.9375 -> 467.6316
$ws.Cells.Item(107, 9).Value = 369.35294  # I107: 235.5 -> 369.35294
$ws.Cells.Item(107, 11).Value = 369.35294  # K107: 235.5 -> 369.35294
$ws.Cells.Item(107, 13).Value = 1550.64706  # M107: 1684.5 -> 1550.64706
$ws.Cells.Item(111, 8).Value = 1852.4  # H111: 2354.3333 -> 1852.4
$ws.Cells.Item(111, 9).Value = 1399.6666  # I111: 2000 -> 1399.6666
$ws.Cells.Item(111, 11).Value = 4198.9998  # K111: 6000 -> 4198.9998
$ws.Cells.Item(111, 13).Value = -1131.9998  # M111: -2933 -> -1131.9998
$ws.Cells.Item(116, 8).Value = 5234.25  # H116: 4277.4287 -> 5234.25
$ws.Cells.Item(116, 9).Value = 5146  # I116: 4073.8333 -> 5146
$ws.Cells.Item(116, 11).Value = 5146  # K116: 4073.8333 -> 5146
$ws.Cells.Item(116, 13).Value = -1704  # M116: -631.8332999999998 -> -1704

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 10975.333  # H28: 9030.25 -> 10975.333
$ws.Cells.Item(28, 9).Value = 10975.333  # I28: 9030.25 -> 10975.333
$ws.Cells.Item(28, 11).Value = 10975.333  # K28: 9030.25 -> 10975.333
$ws.Cells.Item(28, 13).Value = -10783.333  # M28: -8838.25 -> -10783.333
$ws.Cells.Item(61, 8).Value = 12180.833  # H61: 13236.909 -> 12180.833
$ws.Cells.Item(61, 9).Value = 8894.75  # I61: 9370.571 -> 8894.75
$ws.Cells.Item(61, 10).Value = 18753  # J61: 20003 -> 18753
$ws.Cells.Item(61, 11).Value = 8894.75  # K61: 9370.571 -> 8894.75
$ws.Cells.Item(61, 12).Value = 18753  # L61: 20003 -> 18753
$ws.Cells.Item(61, 13).Value = -8682.75  # M61: -9158.571 -> -8682.75
$ws.Cells.Item(61, 14).Value = -19177  # N61: -20427 -> -19177
$ws.Cells.Item(99, 8).Value = 10975.333  # H99: 9030.25 -> 10975.333
$ws.Cells.Item(99, 9).Value = 10975.333  # I99: 9030.25 -> 10975.333
$ws.Cells.Item(99, 11).Value = 10975.333  # K99: 9030.25 -> 10975.333
$ws.Cells.Item(99, 13).Value = -7980.333000000001  # M99: -6035.25 -> -7980.333000000001
$ws.Cells.Item(102, 8).Value = 773.4  # H102: 956 -> 773.4
$ws.Cells.Item(102, 9).Value = 591.75  # I102: 684 -> 591.75
$ws.Cells.Item(102, 11).Value = 591.75  # K102: 684 -> 591.75
$ws.Cells.Item(102, 13).Value = 1030.25  # M102: 938 -> 1030.25
$ws.Cells.Item(136, 8).Value = 12180.833  # H136: 13236.909 -> 12180.833
$ws.Cells.Item(136, 9).Value = 8894.75  # I136: 9370.571 -> 8894.75
$ws.Cells.Item(136, 10).Value = 18753  # J136: 20003 -> 18753
$ws.Cells.Item(136, 11).Value = 26684.25  # K136: 28111.713 -> 26684.25
$ws.Cells.Item(136, 12).Value = 56259  # L136: 60009 -> 56259
$ws.Cells.Item(136, 13).Value = -24134.25  # M136: -25561.713 -> -24134.25
$ws.Cells.Item(136, 14).Value = -61359  # N136: -65109 -> -61359

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 577.5714  # H64: 1002.2 -> 577.5714
$ws.Cells.Item(64, 9).Value = 511.75  # I64: 999 -> 511.75
$ws.Cells.Item(64, 10).Value = 665.3333  # J64: 1004.3333 -> 665.3333
$ws.Cells.Item(64, 11).Value = 511.75  # K64: 999 -> 511.75
$ws.Cells.Item(64, 12).Value = 665.3333  # L64: 1004.3333 -> 665.3333
$ws.Cells.Item(64, 13).Value = -286.75  # M64: -774 -> -286.75
$ws.Cells.Item(64, 14).Value = -1115.3333  # N64: -1454.3333 -> -1115.3333
$ws.Cells.Item(67, 8).Value = 577.5714  # H67: 1002.2 -> 577.5714
$ws.Cells.Item(67, 9).Value = 511.75  # I67: 999 -> 511.75
$ws.Cells.Item(67, 10).Value = 665.3333  # J67: 1004.3333 -> 665.3333
$ws.Cells.Item(67, 11).Value = 511.75  # K67: 999 -> 511.75
$ws.Cells.Item(67, 12).Value = 665.3333  # L67: 1004.3333 -> 665.3333
$ws.Cells.Item(67, 13).Value = 268.25  # M67: -219 -> 268.25
$ws.Cells.Item(67, 14).Value = -2225.3333  # N67: -2564.3333 -> -2225.3333
$ws.Cells.Item(134, 8).Value = 1172.6428  # H134: 917.2 -> 1172.6428
$ws.Cells.Item(134, 9).Value = 1172.6428  # I134: 917.2 -> 1172.6428
$ws.Cells.Item(134, 11).Value = 3517.9284  # K134: 2751.6 -> 3517.9284
$ws.Cells.Item(134, 13).Value = -982.9284000000002  # M134: -216.6000000000004 -> -982.9284000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3035.8  # H58: 3129.7778 -> 3035.8
$ws.Cells.Item(58, 9).Value = 3598.5  # I58: 3880.2 -> 3598.5
$ws.Cells.Item(58, 11).Value = 3598.5  # K58: 3880.2 -> 3598.5
$ws.Cells.Item(58, 13).Value = -3395.5  # M58: -3677.2 -> -3395.5
$ws.Cells.Item(62, 8).Value = 5000  # H62: 3281.5 -> 5000
$ws.Cells.Item(62, 9).Value = 5000  # I62: 3337.8 -> 5000
$ws.Cells.Item(62, 10).Value = 0  # J62: 3000 -> 0
$ws.Cells.Item(62, 11).Value = 5000  # K62: 3337.8 -> 5000
$ws.Cells.Item(62, 12).Value = 0  # L62: 3000 -> 0
$ws.Cells.Item(62, 13).Value = -4376  # M62: -2713.8 -> -4376
$ws.Cells.Item(62, 14).Value = ""  # N62: -4248 -> None
$ws.Cells.Item(65, 8).Value = 5000  # H65: 3281.5 -> 5000
$ws.Cells.Item(65, 9).Value = 5000  # I65: 3337.8 -> 5000
$ws.Cells.Item(65, 10).Value = 0  # J65: 3000 -> 0
$ws.Cells.Item(65, 11).Value = 25000  # K65: 16689 -> 25000
$ws.Cells.Item(65, 12).Value = 0  # L65: 15000 -> 0
$ws.Cells.Item(65, 13).Value = -21880  # M65: -13569 -> -21880
$ws.Cells.Item(65, 14).Value = ""  # N65: -21240 -> None
$ws.Cells.Item(80, 8).Value = 25000  # H80: 15000 -> 25000
$ws.Cells.Item(80, 9).Value = 0  # I80: 10000 -> 0
$ws.Cells.Item(80, 11).Value = 0  # K80: 10000 -> 0
$ws.Cells.Item(80, 13).Value = ""  # M80: -8877 -> None
$ws.Cells.Item(83, 8).Value = 25000  # H83: 15000 -> 25000
$ws.Cells.Item(83, 9).Value = 0  # I83: 10000 -> 0
$ws.Cells.Item(83, 11).Value = 0  # K83: 30000 -> 0
$ws.Cells.Item(83, 13).Value = ""  # M83: -24384 -> None
$ws.Cells.Item(132, 8).Value = 1890.2333  # H132: 1926.1034 -> 1890.2333
$ws.Cells.Item(132, 9).Value = 1996.125  # I132: 2045.9565 -> 1996.125
$ws.Cells.Item(132, 11).Value = 5988.375  # K132: 6137.8695 -> 5988.375
$ws.Cells.Item(132, 13).Value = -3458.375  # M132: -3607.8695 -> -3458.375
$ws.Cells.Item(134, 8).Value = 1892.8422  # H134: 2065 -> 1892.8422
$ws.Cells.Item(134, 9).Value = 1939  # I134: 2173.8333 -> 1939
$ws.Cells.Item(134, 10).Value = 1719.75  # J134: 1629.6666 -> 1719.75
$ws.Cells.Item(134, 11).Value = 5817  # K134: 6521.499899999999 -> 5817
$ws.Cells.Item(134, 12).Value = 5159.25  # L134: 4888.9998 -> 5159.25
$ws.Cells.Item(134, 13).Value = -3282  # M134: -3986.499899999999 -> -3282
$ws.Cells.Item(134, 14).Value = -10229.25  # N134: -9958.9998 -> -10229.25
$ws.Cells.Item(136, 8).Value = 3035.8  # H136: 3129.7778 -> 3035.8
$ws.Cells.Item(136, 9).Value = 3598.5  # I136: 3880.2 -> 3598.5
$ws.Cells.Item(136, 11).Value = 10795.5  # K136: 11640.6 -> 10795.5
$ws.Cells.Item(136, 13).Value = -8245.5  # M136: -9090.599999999999 -> -8245.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 43.42105  # H2: 41.35 -> 43.42105
$ws.Cells.Item(2, 10).Value = 100  # J2: 86 -> 100
$ws.Cells.Item(2, 12).Value = 600  # L2: 516 -> 600
$ws.Cells.Item(2, 14).Value = -826  # N2: -742 -> -826
$ws.Cells.Item(17, 8).Value = 573.8  # H17: 494.83334 -> 573.8
$ws.Cells.Item(38, 8).Value = 83.5  # H38: 39.333332 -> 83.5
$ws.Cells.Item(38, 10).Value = 95.25  # J38: 29 -> 95.25
$ws.Cells.Item(38, 12).Value = 285.75  # L38: 87 -> 285.75
$ws.Cells.Item(38, 14).Value = -979.75  # N38: -781 -> -979.75
$ws.Cells.Item(112, 8).Value = 2454.1667  # H112: 2739 -> 2454.1667
$ws.Cells.Item(112, 9).Value = 2275.3333  # I112: 2913 -> 2275.3333
$ws.Cells.Item(112, 10).Value = 2633  # J112: 2565 -> 2633
$ws.Cells.Item(112, 11).Value = 6825.999899999999  # K112: 8739 -> 6825.999899999999
$ws.Cells.Item(112, 12).Value = 7899  # L112: 7695 -> 7899
$ws.Cells.Item(112, 13).Value = -5717.999899999999  # M112: -7631 -> -5717.999899999999
$ws.Cells.Item(112, 14).Value = -10115  # N112: -9911 -> -10115
$ws.Cells.Item(123, 8).Value = 999  # H123: 0 -> 999
$ws.Cells.Item(123, 9).Value = 999  # I123: 0 -> 999
$ws.Cells.Item(123, 11).Value = 2997  # K123: 0 -> 2997
$ws.Cells.Item(123, 13).Value = -547  # M123: None -> -547

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 25000  # H15: 40000 -> 25000
$ws.Cells.Item(15, 10).Value = 25000  # J15: 40000 -> 25000
$ws.Cells.Item(15, 12).Value = 25000  # L15: 40000 -> 25000
$ws.Cells.Item(15, 14).Value = -25576  # N15: -40576 -> -25576
$ws.Cells.Item(81, 8).Value = 25000  # H81: 40000 -> 25000
$ws.Cells.Item(81, 10).Value = 25000  # J81: 40000 -> 25000
$ws.Cells.Item(81, 12).Value = 25000  # L81: 40000 -> 25000
$ws.Cells.Item(81, 14).Value = -26996  # N81: -41996 -> -26996
$ws.Cells.Item(84, 8).Value = 25000  # H84: 40000 -> 25000
$ws.Cells.Item(84, 10).Value = 25000  # J84: 40000 -> 25000
$ws.Cells.Item(84, 12).Value = 75000  # L84: 120000 -> 75000
$ws.Cells.Item(84, 14).Value = -84984  # N84: -129984 -> -84984
$ws.Cells.Item(113, 8).Value = 3240.1667  # H113: 2492.2222 -> 3240.1667
$ws.Cells.Item(113, 9).Value = 2088.2  # I113: 1748.7142 -> 2088.2
$ws.Cells.Item(113, 10).Value = 9000  # J113: 5094.5 -> 9000
$ws.Cells.Item(113, 11).Value = 2088.2  # K113: 1748.7142 -> 2088.2
$ws.Cells.Item(113, 12).Value = 9000  # L113: 5094.5 -> 9000
$ws.Cells.Item(113, 13).Value = 81.80000000000018  # M113: 421.2858000000001 -> 81.80000000000018
$ws.Cells.Item(113, 14).Value = -13340  # N113: -9434.5 -> -13340

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6316  # H22: 748 -> 6316
$ws.Cells.Item(22, 9).Value = 579.2  # I22: 748 -> 579.2
$ws.Cells.Item(22, 10).Value = 35000  # J22: 0 -> 35000
$ws.Cells.Item(22, 11).Value = 579.2  # K22: 748 -> 579.2
$ws.Cells.Item(22, 12).Value = 35000  # L22: 0 -> 35000
$ws.Cells.Item(22, 13).Value = -284.2  # M22: -453 -> -284.2
$ws.Cells.Item(22, 14).Value = -35590  # N22: None -> -35590
$ws.Cells.Item(27, 8).Value = 6316  # H27: 748 -> 6316
$ws.Cells.Item(27, 9).Value = 579.2  # I27: 748 -> 579.2
$ws.Cells.Item(27, 10).Value = 35000  # J27: 0 -> 35000
$ws.Cells.Item(27, 11).Value = 579.2  # K27: 748 -> 579.2
$ws.Cells.Item(27, 12).Value = 35000  # L27: 0 -> 35000
$ws.Cells.Item(27, 13).Value = -472.2  # M27: -641 -> -472.2
$ws.Cells.Item(27, 14).Value = -35214  # N27: None -> -35214
$ws.Cells.Item(61, 8).Value = 1036  # H61: 1552.7273 -> 1036
$ws.Cells.Item(61, 9).Value = 1036  # I61: 1320 -> 1036
$ws.Cells.Item(61, 10).Value = 0  # J61: 2600 -> 0
$ws.Cells.Item(61, 11).Value = 1036  # K61: 1320 -> 1036
$ws.Cells.Item(61, 12).Value = 0  # L61: 2600 -> 0
$ws.Cells.Item(61, 13).Value = -834  # M61: -1118 -> -834
$ws.Cells.Item(61, 14).Value = ""  # N61: -3004 -> None
$ws.Cells.Item(82, 8).Value = 1279.6  # H82: 1369.0769 -> 1279.6
$ws.Cells.Item(82, 9).Value = 1290.5454  # I82: 1349.8 -> 1290.5454
$ws.Cells.Item(82, 10).Value = 1249.5  # J82: 1433.3334 -> 1249.5
$ws.Cells.Item(82, 11).Value = 1290.5454  # K82: 1349.8 -> 1290.5454
$ws.Cells.Item(82, 12).Value = 1249.5  # L82: 1433.3334 -> 1249.5
$ws.Cells.Item(82, 13).Value = -929.5454  # M82: -988.8 -> -929.5454
$ws.Cells.Item(82, 14).Value = -1971.5  # N82: -2155.3334 -> -1971.5
$ws.Cells.Item(85, 8).Value = 1279.6  # H85: 1369.0769 -> 1279.6
$ws.Cells.Item(85, 9).Value = 1290.5454  # I85: 1349.8 -> 1290.5454
$ws.Cells.Item(85, 10).Value = 1249.5  # J85: 1433.3334 -> 1249.5
$ws.Cells.Item(85, 11).Value = 1290.5454  # K85: 1349.8 -> 1290.5454
$ws.Cells.Item(85, 12).Value = 1249.5  # L85: 1433.3334 -> 1249.5
$ws.Cells.Item(85, 13).Value = -42.54539999999997  # M85: -101.8 -> -42.54539999999997
$ws.Cells.Item(85, 14).Value = -3745.5  # N85: -3929.3334 -> -3745.5
$ws.Cells.Item(92, 8).Value = 0  # H92: 24000 -> 0
$ws.Cells.Item(92, 10).Value = 0  # J92: 24000 -> 0
$ws.Cells.Item(92, 12).Value = 0  # L92: 24000 -> 0
$ws.Cells.Item(92, 14).Value = ""  # N92: -28992 -> None
$ws.Cells.Item(96, 8).Value = 0  # H96: 33999 -> 0
$ws.Cells.Item(96, 10).Value = 0  # J96: 33999 -> 0
$ws.Cells.Item(96, 12).Value = 0  # L96: 33999 -> 0
$ws.Cells.Item(96, 14).Value = ""  # N96: -39491 -> None
$ws.Cells.Item(99, 8).Value = 15272.333  # H99: 15703.75 -> 15272.333
$ws.Cells.Item(99, 9).Value = 15272.333  # I99: 15703.75 -> 15272.333
$ws.Cells.Item(99, 11).Value = 15272.333  # K99: 15703.75 -> 15272.333
$ws.Cells.Item(99, 13).Value = -12277.333  # M99: -12708.75 -> -12277.333
$ws.Cells.Item(101, 8).Value = 21614.834  # H101: 19003.715 -> 21614.834
$ws.Cells.Item(101, 10).Value = 21614.834  # J101: 19003.715 -> 21614.834
$ws.Cells.Item(101, 12).Value = 21614.834  # L101: 19003.715 -> 21614.834
$ws.Cells.Item(101, 14).Value = -28104.834  # N101: -25493.715 -> -28104.834
$ws.Cells.Item(113, 8).Value = 1036  # H113: 1552.7273 -> 1036
$ws.Cells.Item(113, 9).Value = 1036  # I113: 1320 -> 1036
$ws.Cells.Item(113, 10).Value = 0  # J113: 2600 -> 0
$ws.Cells.Item(113, 11).Value = 1036  # K113: 1320 -> 1036
$ws.Cells.Item(113, 12).Value = 0  # L113: 2600 -> 0
$ws.Cells.Item(113, 13).Value = 1134  # M113: 850 -> 1134
$ws.Cells.Item(113, 14).Value = ""  # N113: -6940 -> None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 0  # H5: 600000000 -> 0
$ws.Cells.Item(5, 9).Value = 0  # I5: 600000000 -> 0
$ws.Cells.Item(5, 11).Value = 0  # K5: 600000000 -> 0
$ws.Cells.Item(5, 13).Value = ""  # M5: -599999888 -> None
$ws.Cells.Item(122, 8).Value = 3394.111  # H122: 3690.5 -> 3394.111
$ws.Cells.Item(122, 9).Value = 3292.4285  # I122: 3646.2856 -> 3292.4285
$ws.Cells.Item(122, 10).Value = 3750  # J122: 4000 -> 3750
$ws.Cells.Item(122, 11).Value = 9877.2855  # K122: 10938.8568 -> 9877.2855
$ws.Cells.Item(122, 12).Value = 11250  # L122: 12000 -> 11250
$ws.Cells.Item(122, 13).Value = -7427.2855  # M122: -8488.856800000001 -> -7427.2855
$ws.Cells.Item(122, 14).Value = -16150  # N122: -16900 -> -16150
$ws.Cells.Item(132, 8).Value = 4391.6665  # H132: 3074.875 -> 4391.6665
$ws.Cells.Item(132, 9).Value = 5120  # I132: 3414.1428 -> 5120
$ws.Cells.Item(132, 10).Value = 750  # J132: 700 -> 750
$ws.Cells.Item(132, 11).Value = 15360  # K132: 10242.4284 -> 15360
$ws.Cells.Item(132, 12).Value = 2250  # L132: 2100 -> 2250
$ws.Cells.Item(132, 13).Value = -12830  # M132: -7712.428400000001 -> -12830
$ws.Cells.Item(132, 14).Value = -7310  # N132: -7160 -> -7310
